$d = $word.ActiveDocument

$d.Content.Find.Execute("60×51=3060", $true, $false, $false, $false, $false, $true, 1, $false, "36×59=2124", 2) | Out-Null
$d.Content.Find.Execute("57×66=3762", $true, $false, $false, $false, $false, $true, 1, $false, "69×18=1242", 2) | Out-Null
$d.Content.Find.Execute("68×35=2380", $true, $false, $false, $false, $false, $true, 1, $false, "52×17=884", 2) | Out-Null
$d.Content.Find.Execute("46×88=4048", $true, $false, $false, $false, $false, $true, 1, $false, "69×57=3933", 2) | Out-Null
$d.Content.Find.Execute("65×19=1235", $true, $false, $false, $false, $false, $true, 1, $false, "30×14=420", 2) | Out-Null
$d.Content.Find.Execute("54×99=5346", $true, $false, $false, $false, $false, $true, 1, $false, "95×18=1710", 2) | Out-Null
$d.Content.Find.Execute("10×58=580", $true, $false, $false, $false, $false, $true, 1, $false, "33×78=2574", 2) | Out-Null
$d.Content.Find.Execute("57×78=4446", $true, $false, $false, $false, $false, $true, 1, $false, "91×62=5642", 2) | Out-Null
$d.Content.Find.Execute("42×81=3402", $true, $false, $false, $false, $false, $true, 1, $false, "45×71=3195", 2) | Out-Null
$d.Content.Find.Execute("14×90=1260", $true, $false, $false, $false, $false, $true, 1, $false, "98×75=7350", 2) | Out-Null
$d.Content.Find.Execute("37×97=3589", $true, $false, $false, $false, $false, $true, 1, $false, "47×18=846", 2) | Out-Null
$d.Content.Find.Execute("59×46=2714", $true, $false, $false, $false, $false, $true, 1, $false, "29×85=2465", 2) | Out-Null
$d.Content.Find.Execute("84×99=8316", $true, $false, $false, $false, $false, $true, 1, $false, "41×21=861", 2) | Out-Null
$d.Content.Find.Execute("87×79=6873", $true, $false, $false, $false, $false, $true, 1, $false, "47×46=2162", 2) | Out-Null
$d.Content.Find.Execute("45×41=1845", $true, $false, $false, $false, $false, $true, 1, $false, "64×61=3904", 2) | Out-Null
$d.Content.Find.Execute("13×79=1027", $true, $false, $false, $false, $false, $true, 1, $false, "43×50=2150", 2) | Out-Null
$d.Content.Find.Execute("49×72=3528", $true, $false, $false, $false, $false, $true, 1, $false, "36×69=2484", 2) | Out-Null
$d.Content.Find.Execute("52×39=2028", $true, $false, $false, $false, $false, $true, 1, $false, "37×50=1850", 2) | Out-Null
$d.Content.Find.Execute("57×63=3591", $true, $false, $false, $false, $false, $true, 1, $false, "83×100=8300", 2) | Out-Null
$d.Content.Find.Execute("91×59=5369", $true, $false, $false, $false, $false, $true, 1, $false, "29×25=725", 2) | Out-Null
$d.Content.Find.Execute("16×54=864", $true, $false, $false, $false, $false, $true, 1, $false, "75×26=1950", 2) | Out-Null
$d.Content.Find.Execute("93×35=3255", $true, $false, $false, $false, $false, $true, 1, $false, "96×17=1632", 2) | Out-Null
$d.Content.Find.Execute("50×31=1550", $true, $false, $false, $false, $false, $true, 1, $false, "11×17=187", 2) | Out-Null
$d.Content.Find.Execute("37×86=3182", $true, $false, $false, $false, $false, $true, 1, $false, "33×10=330", 2) | Out-Null
$d.Content.Find.Execute("97×10=970", $true, $false, $false, $false, $false, $true, 1, $false, "38×10=380", 2) | Out-Null
$d.Content.Find.Execute("81×81=6561", $true, $false, $false, $false, $false, $true, 1, $false, "37×19=703", 2) | Out-Null
$d.Content.Find.Execute("61×84=5124", $true, $false, $false, $false, $false, $true, 1, $false, "83×33=2739", 2) | Out-Null
$d.Content.Find.Execute("27×96=2592", $true, $false, $false, $false, $false, $true, 1, $false, "57×49=2793", 2) | Out-Null
$d.Content.Find.Execute("59×11=649", $true, $false, $false, $false, $false, $true, 1, $false, "37×87=3219", 2) | Out-Null
$d.Content.Find.Execute("53×93=4929", $true, $false, $false, $false, $false, $true, 1, $false, "14×55=770", 2) | Out-Null
$d.Content.Find.Execute("64×94=6016", $true, $false, $false, $false, $false, $true, 1, $false, "25×16=400", 2) | Out-Null
$d.Content.Find.Execute("32×76=2432", $true, $false, $false, $false, $false, $true, 1, $false, "53×73=3869", 2) | Out-Null
$d.Content.Find.Execute("53×84=4452", $true, $false, $false, $false, $false, $true, 1, $false, "43×19=817", 2) | Out-Null
$d.Content.Find.Execute("82×94=7708", $true, $false, $false, $false, $false, $true, 1, $false, "64×83=5312", 2) | Out-Null
$d.Content.Find.Execute("84×28=2352", $true, $false, $false, $false, $false, $true, 1, $false, "13×86=1118", 2) | Out-Null
$d.Content.Find.Execute("92×51=4692", $true, $false, $false, $false, $false, $true, 1, $false, "66×41=2706", 2) | Out-Null
$d.Content.Find.Execute("90×53=4770", $true, $false, $false, $false, $false, $true, 1, $false, "32×21=672", 2) | Out-Null
$d.Content.Find.Execute("68×94=6392", $true, $false, $false, $false, $false, $true, 1, $false, "93×47=4371", 2) | Out-Null
$d.Content.Find.Execute("49×76=3724", $true, $false, $false, $false, $false, $true, 1, $false, "80×92=7360", 2) | Out-Null
$d.Content.Find.Execute("66×17=1122", $true, $false, $false, $false, $false, $true, 1, $false, "86×84=7224", 2) | Out-Null
$d.Content.Find.Execute("27×21=567", $true, $false, $false, $false, $false, $true, 1, $false, "51×14=714", 2) | Out-Null
$d.Content.Find.Execute("91×73=6643", $true, $false, $false, $false, $false, $true, 1, $false, "18×23=414", 2) | Out-Null
$d.Content.Find.Execute("25×81=2025", $true, $false, $false, $false, $false, $true, 1, $false, "70×21=1470", 2) | Out-Null
$d.Content.Find.Execute("40×18=720", $true, $false, $false, $false, $false, $true, 1, $false, "67×55=3685", 2) | Out-Null
$d.Content.Find.Execute("19×82=1558", $true, $false, $false, $false, $false, $true, 1, $false, "69×47=3243", 2) | Out-Null
$d.Content.Find.Execute("21×93=1953", $true, $false, $false, $false, $false, $true, 1, $false, "98×74=7252", 2) | Out-Null
$d.Content.Find.Execute("50×98=4900", $true, $false, $false, $false, $false, $true, 1, $false, "31×28=868", 2) | Out-Null
$d.Content.Find.Execute("94×53=4982", $true, $false, $false, $false, $false, $true, 1, $false, "22×94=2068", 2) | Out-Null
$d.Content.Find.Execute("89×14=1246", $true, $false, $false, $false, $false, $true, 1, $false, "32×80=2560", 2) | Out-Null
$d.Content.Find.Execute("98×30=2940", $true, $false, $false, $false, $false, $true, 1, $false, "30×10=300", 2) | Out-Null
$d.Content.Find.Execute("35×91=3185", $true, $false, $false, $false, $false, $true, 1, $false, "21×82=1722", 2) | Out-Null
$d.Content.Find.Execute("17×14=238", $true, $false, $false, $false, $false, $true, 1, $false, "98×77=7546", 2) | Out-Null
$d.Content.Find.Execute("63×56=3528", $true, $false, $false, $false, $false, $true, 1, $false, "20×21=420", 2) | Out-Null
$d.Content.Find.Execute("70×41=2870", $true, $false, $false, $false, $false, $true, 1, $false, "63×17=1071", 2) | Out-Null
$d.Content.Find.Execute("11×79=869", $true, $false, $false, $false, $false, $true, 1, $false, "84×79=6636", 2) | Out-Null
$d.Content.Find.Execute("57×71=4047", $true, $false, $false, $false, $false, $true, 1, $false, "33×46=1518", 2) | Out-Null
$d.Content.Find.Execute("25×33=825", $true, $false, $false, $false, $false, $true, 1, $false, "65×97=6305", 2) | Out-Null
$d.Content.Find.Execute("30×13=390", $true, $false, $false, $false, $false, $true, 1, $false, "97×40=3880", 2) | Out-Null
$d.Content.Find.Execute("55×43=2365", $true, $false, $false, $false, $false, $true, 1, $false, "72×19=1368", 2) | Out-Null
$d.Content.Find.Execute("38×43=1634", $true, $false, $false, $false, $false, $true, 1, $false, "60×45=2700", 2) | Out-Null
$d.Content.Find.Execute("45×23=1035", $true, $false, $false, $false, $false, $true, 1, $false, "32×71=2272", 2) | Out-Null
$d.Content.Find.Execute("92×50=4600", $true, $false, $false, $false, $false, $true, 1, $false, "34×14=476", 2) | Out-Null
$d.Content.Find.Execute("55×55=3025", $true, $false, $false, $false, $false, $true, 1, $false, "64×93=5952", 2) | Out-Null
$d.Content.Find.Execute("61×55=3355", $true, $false, $false, $false, $false, $true, 1, $false, "93×11=1023", 2) | Out-Null
$d.Content.Find.Execute("99×68=6732", $true, $false, $false, $false, $false, $true, 1, $false, "80×51=4080", 2) | Out-Null
$d.Content.Find.Execute("46×70=3220", $true, $false, $false, $false, $false, $true, 1, $false, "97×87=8439", 2) | Out-Null
$d.Content.Find.Execute("57×27=1539", $true, $false, $false, $false, $false, $true, 1, $false, "36×97=3492", 2) | Out-Null
$d.Content.Find.Execute("84×63=5292", $true, $false, $false, $false, $false, $true, 1, $false, "21×30=630", 2) | Out-Null
$d.Content.Find.Execute("100×44=4400", $true, $false, $false, $false, $false, $true, 1, $false, "70×35=2450", 2) | Out-Null
$d.Content.Find.Execute("25×88=2200", $true, $false, $false, $false, $false, $true, 1, $false, "37×49=1813", 2) | Out-Null
$d.Content.Find.Execute("22×15=330", $true, $false, $false, $false, $false, $true, 1, $false, "56×34=1904", 2) | Out-Null
$d.Content.Find.Execute("41×47=1927", $true, $false, $false, $false, $false, $true, 1, $false, "56×11=616", 2) | Out-Null
$d.Content.Find.Execute("11×75=825", $true, $false, $false, $false, $false, $true, 1, $false, "39×13=507", 2) | Out-Null
$d.Content.Find.Execute("61×95=5795", $true, $false, $false, $false, $false, $true, 1, $false, "46×73=3358", 2) | Out-Null
$d.Content.Find.Execute("42×90=3780", $true, $false, $false, $false, $false, $true, 1, $false, "24×35=840", 2) | Out-Null
$d.Content.Find.Execute("66×99=6534", $true, $false, $false, $false, $false, $true, 1, $false, "100×41=4100", 2) | Out-Null
$d.Content.Find.Execute("13×59=767", $true, $false, $false, $false, $false, $true, 1, $false, "42×32=1344", 2) | Out-Null
$d.Content.Find.Execute("74×39=2886", $true, $false, $false, $false, $false, $true, 1, $false, "46×67=3082", 2) | Out-Null
$d.Content.Find.Execute("93×25=2325", $true, $false, $false, $false, $false, $true, 1, $false, "33×46=1518", 2) | Out-Null
$d.Content.Find.Execute("58×46=2668", $true, $false, $false, $false, $false, $true, 1, $false, "70×66=4620", 2) | Out-Null
$d.Content.Find.Execute("78×72=5616", $true, $false, $false, $false, $false, $true, 1, $false, "55×64=3520", 2) | Out-Null
$d.Content.Find.Execute("66×98=6468", $true, $false, $false, $false, $false, $true, 1, $false, "66×11=726", 2) | Out-Null
$d.Content.Find.Execute("80×44=3520", $true, $false, $false, $false, $false, $true, 1, $false, "61×90=5490", 2) | Out-Null
$d.Content.Find.Execute("13×69=897", $true, $false, $false, $false, $false, $true, 1, $false, "44×94=4136", 2) | Out-Null
$d.Content.Find.Execute("92×42=3864", $true, $false, $false, $false, $false, $true, 1, $false, "27×39=1053", 2) | Out-Null
$d.Content.Find.Execute("77×54=4158", $true, $false, $false, $false, $false, $true, 1, $false, "91×25=2275", 2) | Out-Null
$d.Content.Find.Execute("48×63=3024", $true, $false, $false, $false, $false, $true, 1, $false, "65×44=2860", 2) | Out-Null
$d.Content.Find.Execute("71×95=6745", $true, $false, $false, $false, $false, $true, 1, $false, "34×11=374", 2) | Out-Null
$d.Content.Find.Execute("78×86=6708", $true, $false, $false, $false, $false, $true, 1, $false, "37×40=1480", 2) | Out-Null
$d.Content.Find.Execute("58×24=1392", $true, $false, $false, $false, $false, $true, 1, $false, "91×96=8736", 2) | Out-Null
$d.Content.Find.Execute("47×48=2256", $true, $false, $false, $false, $false, $true, 1, $false, "39×62=2418", 2) | Out-Null
$d.Content.Find.Execute("26×53=1378", $true, $false, $false, $false, $false, $true, 1, $false, "47×58=2726", 2) | Out-Null
$d.Content.Find.Execute("32×66=2112", $true, $false, $false, $false, $false, $true, 1, $false, "94×49=4606", 2) | Out-Null
$d.Content.Find.Execute("95×32=3040", $true, $false, $false, $false, $false, $true, 1, $false, "26×90=2340", 2) | Out-Null
$d.Content.Find.Execute("68×91=6188", $true, $false, $false, $false, $false, $true, 1, $false, "42×11=462", 2) | Out-Null
$d.Content.Find.Execute("76×95=7220", $true, $false, $false, $false, $false, $true, 1, $false, "78×52=4056", 2) | Out-Null
$d.Content.Find.Execute("25×85=2125", $true, $false, $false, $false, $false, $true, 1, $false, "100×99=9900", 2) | Out-Null
$d.Content.Find.Execute("13×73=949", $true, $false, $false, $false, $false, $true, 1, $false, "11×42=462", 2) | Out-Null
$d.Content.Find.Execute("24×62=1488", $true, $false, $false, $false, $false, $true, 1, $false, "85×49=4165", 2) | Out-Null
$d.Content.Find.Execute("35×38=1330", $true, $false, $false, $false, $false, $true, 1, $false, "75×96=7200", 2) | Out-Null
